$wb = $excel.ActiveWorkbook

# The repayment strategy scenario on the ProductLoanInput sheet is being
# changed from "RBI (India)" to "Overdue/Due Fee/Int,Principal" as part of
# adding periodic & upfront related scenarios.
$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsInput.Activate()

$wsInput.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Leave the selection on the cell that was just updated, and make this the
# active/selected sheet of the workbook.
$wsInput.Range("B17").Select()
